$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 00:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 584862
$ws.Range("C4").Value = 24562
$ws.Range("D4").Value = 36205
$ws.Range("E4").Value = 525102
$ws.Range("G4").Value = 1450
$ws.Range("H4").Value = 23555

# Row 16 - Canada
$ws.Range("D16").Value = 7756
$ws.Range("E16").Value = 17029

# Row 32 - Noruega
$ws.Range("B32").Value = 6603
$ws.Range("C32").Value = 78
$ws.Range("E32").Value = 6437

# Row 52 - Colombia
$ws.Range("B52").Value = 2852
$ws.Range("C52").Value = 76
$ws.Range("D52").Value = 319
$ws.Range("E52").Value = 2421
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 112

# Row 140 - Jamaica
$ws.Range("B140").Value = 73
$ws.Range("C140").Value = 1
$ws.Range("D140").Value = 19
$ws.Range("E140").Value = 50

# Row 141 - Barbados
$ws.Range("B141").Value = 72
$ws.Range("C141").Value = 1
$ws.Range("D141").Value = 13
$ws.Range("E141").Value = 55
